$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 2429.1428
$ws.Range("I40").Value = 1999.6666
$ws.Range("K40").Value = 1999.6666
$ws.Range("M40").Value = -1824.6666

$ws.Range("H62").Value = 1150
$ws.Range("I62").Value = 900
$ws.Range("K62").Value = 900
$ws.Range("M62").Value = -276

$ws.Range("H65").Value = 1150
$ws.Range("I65").Value = 900
$ws.Range("K65").Value = 4500
$ws.Range("M65").Value = -1380

$ws.Range("H74").Value = 5000
$ws.Range("I74").Value = 5000
$ws.Range("K74").Value = 5000
$ws.Range("M74").Value = -4064

$ws.Range("H77").Value = 5000
$ws.Range("I77").Value = 5000
$ws.Range("K77").Value = 25000
$ws.Range("M77").Value = -20320

$ws.Range("H116").Value = 3050
$ws.Range("J116").Value = 2950
$ws.Range("L116").Value = 2950
$ws.Range("N116").Value = -9834

$ws.Range("H137").Value = 1234.8948
$ws.Range("I137").Value = 791.94116
$ws.Range("K137").Value = 2375.82348
$ws.Range("M137").Value = 174.17652

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H33").Value = 5675.3335
$ws.Range("I33").Value = 2513
$ws.Range("K33").Value = 2513
$ws.Range("M33").Value = -2184

$ws.Range("H34").Value = 0
$ws.Range("J34").Value = 0
$ws.Range("L34").Value = 0
$ws.Range("N34").ClearContents()

$ws.Range("H38").Value = 2006
$ws.Range("I38").Value = 2006
$ws.Range("K38").Value = 2006
$ws.Range("M38").Value = -1539

$ws.Range("H59").Value = 6000
$ws.Range("J59").Value = 6000
$ws.Range("L59").Value = 6000
$ws.Range("N59").Value = -7608

$ws.Range("H61").Value = 15002
$ws.Range("I61").Value = 15002
$ws.Range("K61").Value = 15002
$ws.Range("M61").Value = -14790

$ws.Range("H74").Value = 2145.8333
$ws.Range("J74").Value = 2998.3333
$ws.Range("L74").Value = 2998.3333
$ws.Range("N74").Value = -4746.3333

$ws.Range("H77").Value = 2145.8333
$ws.Range("J77").Value = 2998.3333
$ws.Range("L77").Value = 14991.6665
$ws.Range("N77").Value = -23727.6665

$ws.Range("H132").Value = 5897.2
$ws.Range("I132").Value = 4871.625
$ws.Range("K132").Value = 14614.875
$ws.Range("M132").Value = -12084.875

$ws.Range("H136").Value = 15002
$ws.Range("I136").Value = 15002
$ws.Range("K136").Value = 45006
$ws.Range("M136").Value = -42456

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H37").Value = 1300
$ws.Range("J37").Value = 1300
$ws.Range("L37").Value = 1300
$ws.Range("N37").Value = -1574

$ws.Range("H108").Value = 60000
$ws.Range("J108").Value = 60000
$ws.Range("L108").Value = 60000
$ws.Range("N108").Value = -67680

$ws.Range("H134").Value = 3134.8462
$ws.Range("I134").Value = 2477.5454
$ws.Range("K134").Value = 7432.6362
$ws.Range("M134").Value = -4897.6362

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2662.1562
$ws.Range("I31").Value = 1591.7084
$ws.Range("J31").Value = 5873.5
$ws.Range("K31").Value = 1591.7084
$ws.Range("L31").Value = 5873.5
$ws.Range("M31").Value = -1296.7084
$ws.Range("N31").Value = -6463.5

$ws.Range("H34").Value = 2662.1562
$ws.Range("I34").Value = 1591.7084
$ws.Range("J34").Value = 5873.5
$ws.Range("K34").Value = 1591.7084
$ws.Range("L34").Value = 5873.5
$ws.Range("M34").Value = -1389.7084
$ws.Range("N34").Value = -6277.5

$ws.Range("H58").Value = 3286.5
$ws.Range("I58").Value = 2236.75
$ws.Range("K58").Value = 2236.75
$ws.Range("M58").Value = -2033.75

$ws.Range("H132").Value = 2091.95
$ws.Range("I132").Value = 1961.1177
$ws.Range("K132").Value = 5883.3531
$ws.Range("M132").Value = -3353.3531

$ws.Range("H136").Value = 3286.5
$ws.Range("I136").Value = 2236.75
$ws.Range("K136").Value = 6710.25
$ws.Range("M136").Value = -4160.25

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 515.8
$ws.Range("J5").Value = 998
$ws.Range("L5").Value = 2994
$ws.Range("N5").Value = -3218

$ws.Range("H75").Value = 1406.6666
$ws.Range("I75").Value = 1402.5
$ws.Range("J75").Value = 1415
$ws.Range("K75").Value = 4207.5
$ws.Range("L75").Value = 4245
$ws.Range("M75").Value = -3209.5
$ws.Range("N75").Value = -6241

$ws.Range("H78").Value = 1406.6666
$ws.Range("I78").Value = 1402.5
$ws.Range("J78").Value = 1415
$ws.Range("K78").Value = 12622.5
$ws.Range("L78").Value = 12735
$ws.Range("M78").Value = -7630.5
$ws.Range("N78").Value = -22719

$ws.Range("H106").Value = 9000
$ws.Range("J106").Value = 9000
$ws.Range("L106").Value = 27000
$ws.Range("N106").Value = -28892

$ws.Range("H132").Value = 3347.923
$ws.Range("I132").Value = 3069.7778
$ws.Range("J132").Value = 3973.75
$ws.Range("K132").Value = 27628.0002
$ws.Range("L132").Value = 35763.75
$ws.Range("M132").Value = -25098.0002
$ws.Range("N132").Value = -40823.75

$ws.Range("H135").Value = 515.8
$ws.Range("J135").Value = 998
$ws.Range("L135").Value = 8982
$ws.Range("N135").Value = -14052

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H49").Value = 70030
$ws.Range("J49").Value = 70030
$ws.Range("L49").Value = 70030
$ws.Range("N49").Value = -70398

$ws.Range("H122").Value = 66213.625
$ws.Range("I122").Value = 2998.2727
$ws.Range("J122").Value = 205287.4
$ws.Range("K122").Value = 8994.8181
$ws.Range("L122").Value = 615862.2
$ws.Range("M122").Value = -6544.8181
$ws.Range("N122").Value = -620762.2

$ws.Range("H132").Value = 4645.3076
$ws.Range("I132").Value = 4645.3076
$ws.Range("K132").Value = 13935.9228
$ws.Range("M132").Value = -11405.9228

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H4").Value = 7572.6665
$ws.Range("J4").Value = 15000
$ws.Range("L4").Value = 15000
$ws.Range("N4").Value = -15226

$ws.Range("H22").Value = 8463.454
$ws.Range("I22").Value = 3381.3635
$ws.Range("J22").Value = 13545.546
$ws.Range("K22").Value = 3381.3635
$ws.Range("L22").Value = 13545.546
$ws.Range("M22").Value = -3086.3635
$ws.Range("N22").Value = -14135.546

$ws.Range("H27").Value = 8463.454
$ws.Range("I27").Value = 3381.3635
$ws.Range("J27").Value = 13545.546
$ws.Range("K27").Value = 3381.3635
$ws.Range("L27").Value = 13545.546
$ws.Range("M27").Value = -3274.3635
$ws.Range("N27").Value = -13759.546

$ws.Range("H28").Value = 7572.6665
$ws.Range("J28").Value = 15000
$ws.Range("L28").Value = 15000
$ws.Range("N28").Value = -15464

$ws.Range("H36").Value = 100000
$ws.Range("J36").Value = 100000
$ws.Range("L36").Value = 100000
$ws.Range("N36").Value = -101124

$ws.Range("H37").Value = 7572.6665
$ws.Range("J37").Value = 15000
$ws.Range("L37").Value = 15000
$ws.Range("N37").Value = -15214

$ws.Range("H43").Value = 2865714.2
$ws.Range("I43").Value = 9999
$ws.Range("J43").Value = 3341666.8
$ws.Range("K43").Value = 9999
$ws.Range("L43").Value = 3341666.8
$ws.Range("M43").Value = -9806
$ws.Range("N43").Value = -3342052.8

$ws.Range("H51").Value = 20077
$ws.Range("I51").Value = 20077
$ws.Range("K51").Value = 20077
$ws.Range("M51").Value = -19599

$ws.Range("H58").Value = 450
$ws.Range("I58").Value = 450
$ws.Range("K58").Value = 450
$ws.Range("M58").Value = -190

$ws.Range("H68").Value = 7500
$ws.Range("J68").Value = 7500
$ws.Range("L68").Value = 7500
$ws.Range("N68").Value = -8998

$ws.Range("H71").Value = 7500
$ws.Range("J71").Value = 7500
$ws.Range("L71").Value = 37500
$ws.Range("N71").Value = -44988

$ws.Range("H122").Value = 6976.2856
$ws.Range("I122").Value = 6935.8
$ws.Range("J122").Value = 6998.778
$ws.Range("K122").Value = 20807.4
$ws.Range("L122").Value = 20996.334
$ws.Range("M122").Value = -18357.4
$ws.Range("N122").Value = -25896.334
